$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ B=..; C=..; D=..; E=.. } giving the final values that
# need to land in that row after the update. Only columns that actually
# change are listed; rows 36-41 are a full re-shuffle of coins, the rest
# are simple price/volume refreshes.
$updates = @{
    2  = @{ D = "34.511.71"; E = "  +0.46%  " }
    3  = @{ D = "1.801.49";  E = "  +0.21%  " }
    4  = @{ E = "  +0.34%  " }
    5  = @{ D = "223.95";    E = "  -0.54%  " }
    6  = @{ D = "0.599";     E = "  -0.02%  " }
    7  = @{ E = "  +0.28%  " }
    8  = @{ D = "41.17";     E = "  +13.45%  " }
    9  = @{ D = "0.291";     E = "  -0.13%  " }
    10 = @{ D = "0.0665";    E = "  -1.39%  " }
    11 = @{ D = "0.0997";    E = "  +3.50%  " }
    12 = @{ D = "2.062.43";  E = "  +0.28%  " }
    13 = @{ D = "1.803.10";  E = "  +0.06%  " }
    14 = @{ D = "10.92";     E = "  -2.78%  " }
    15 = @{ D = "34.506.13"; E = "  +0.57%  " }
    16 = @{ E = "  -0.30%  " }
    17 = @{ E = "  -0.38%  " }
    18 = @{ D = "67.15";     E = "  -2.16%  " }
    19 = @{ D = "240.46";    E = "  -2.16%  " }
    20 = @{ E = "  -0.62%  " }
    21 = @{ E = "  -1.73%  " }
    22 = @{ E = "  +0.23%  " }
    23 = @{ D = "4.27";      E = "  +4.92%  " }
    24 = @{ D = "2.15";      E = "  -2.56%  " }
    25 = @{ D = "172.03";    E = "  +0.95%  " }
    26 = @{ D = "7.67";      E = "  -2.50%  " }
    27 = @{ D = "17.37";     E = "  +0.28%  " }
    28 = @{ E = "  +0.35%  " }
    29 = @{ E = "  +0.43%  " }
    30 = @{ D = "3.79";      E = "  +0.33%  " }
    31 = @{ E = "  -0.44%  " }
    32 = @{ D = "3.86";      E = "  -1.04%  " }
    33 = @{ E = "  -0.13%  " }
    34 = @{ E = "  +1.47%  " }
    35 = @{ D = "0.648";     E = "  +0.18%  " }

    # Rows 36-41 are reshuffled: Aave moves to the top of this block,
    # followed by Maker, TrustWalletToken, InjectiveProtocol, VeChain,
    # RenderToken (each also gets refreshed price/volume figures).
    36 = @{ B = "Aave";              C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave";                    D = "87.27";   E = "  +8.08%  " }
    37 = @{ B = "Maker";             C = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr";                D = "1.323.56"; E = "  -2.76%  " }
    38 = @{ B = "TrustWalletToken";  C = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt";         D = "1.06";    E = "  +0.84%  " }
    39 = @{ B = "InjectiveProtocol"; C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj";        D = "14.77";   E = "  +11.88%  " }
    40 = @{ B = "VeChain";           C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet";              D = "0.0187";  E = "  +0.79%  " }
    41 = @{ B = "RenderToken";       C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr";         D = "2.34";    E = "  -0.75%  " }

    42 = @{ D = "1.22";      E = "  +5.41%  " }
    43 = @{ E = "  +0.43%  " }
    44 = @{ E = "  +0.15%  " }
    45 = @{ D = "0.935";     E = "  -0.18%  " }
    46 = @{ D = "0.0518";    E = "  +4.45%  " }
    47 = @{ D = "1.963.40";  E = "  +0.31%  " }
    48 = @{ D = "5.79";      E = "  +0.68%  " }
    49 = @{ E = "  +0.19%  " }
    50 = @{ D = "100.50";    E = "  -1.30%  " }
    51 = @{ D = "0.0608";    E = "  +0.93%  " }
}

foreach ($rowNum in $updates.Keys) {
    $rowVals = $updates[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $cell = $ws.Range("$col$rowNum")
        if ($col -eq "D") {
            # The "Price" column holds values like "34.511.71" or "223.95"
            # that must stay as plain text (that's how the sheet already
            # stores them, as inlineStr). Excel auto-coerces
            # plain-decimal-looking strings (e.g. "223.95") to numbers on
            # assignment, so force Text format first on just this cell to
            # preserve the original string formatting/precision without
            # touching the format of cells whose price didn't change.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowVals[$col]
    }
}
